$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole data range to Text format first so Excel does not
# reinterpret numeric-looking strings (e.g. "63.418.03", "6.00") as
# numbers/dates when we assign them below.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.418.03"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.418.28"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "570.06"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").Value = "156.03"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "3.419.61"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").Value = "7.42"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "4.006.56"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").Value = "0.0000190"
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("D16").Value = "27.07"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "63.618.01"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "3.399.41"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("D20").Value = "14.14"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").Value = "386.63"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "8.20"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "72.10"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.536"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").Value = "  +21.72%  "
$ws.Range("D27").Value = "9.48"
$ws.Range("E27").Value = "  +7.30%  "
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "6.00"
$ws.Range("E30").Value = "  +6.78%  "
$ws.Range("D31").Value = "1.36"
$ws.Range("E31").Value = "  +4.45%  "
$ws.Range("D32").Value = "2.00"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "23.26"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "6.43"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "6.83"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").Value = "159.12"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "1.46"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").Value = "0.0764"
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("D40").Value = "1.83"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("D41").Value = "2.896.15"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("D42").Value = "26.83"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "0.0318"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").Value = "4.38"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").Value = "0.763"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").Value = "40.97"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "23.58"
$ws.Range("E47").Value = "  +6.74%  "
$ws.Range("E48").Value = "  +3.49%  "
$ws.Range("D49").Value = "2.19"
$ws.Range("E49").Value = "  +20.51%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "6.44"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "0.841"
$ws.Range("E51").Value = "  +3.93%  "

# Restore default (General) formatting so the saved styles exactly
# match the original workbook (only cell contents were changed).
$dataRange.ClearFormats()
